$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data table for Trial 1..16 (Trial, Question, ConditionType, ITI)
$data = @(
    @(1,10,1,8),
    @(2,3,1,7),
    @(3,35,4,8),
    @(4,29,4,7),
    @(5,30,3,8),
    @(6,6,3,6),
    @(7,18,4,8),
    @(8,31,1,6),
    @(9,8,2,6),
    @(10,32,2,7),
    @(11,9,3,7),
    @(12,5,4,6),
    @(13,7,2,9),
    @(14,27,1,6),
    @(15,23,3,6),
    @(16,22,2,7)
)

# Add the new header for column D
$ws.Cells.Item(1, 4).Value = "ITI"

# Write out the data rows (rows 2..17), updating existing columns and
# populating the new column D
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Remove the now-unused trailing rows (previously Trial 17, 18, 19)
$ws.Range("A18:D20").Clear()

$ws.Range("G11").Select()
